# Update "想去人数" (F column) figures across all four sheets to match the
# newly generated gh-pages snapshot (commit 456a3b4). Row numbers below are
# the actual worksheet rows; each tuple is (row, newValue).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$sheet1Edits = @(
    @(3, 1138),
    @(4, 1217),
    @(6, 163),
    @(7, 529),
    @(8, 295),
    @(9, 50),
    @(10, 1250),
    @(12, 3292),
    @(14, 237),
    @(15, 456),
    @(16, 10),
    @(18, 2),
    @(19, 309),
    @(20, 593),
    @(21, 264),
    @(22, 250),
    @(23, 339),
    @(25, 34),
    @(26, 651),
    @(27, 194),
    @(28, 92),
    @(29, 506),
    @(30, 70),
    @(32, 597)
)
foreach ($edit in $sheet1Edits) {
    $row = $edit[0]
    $val = $edit[1]
    $ws1.Cells.Item($row, 6).Value = $val
}

# Row 11 on sheet 1 also picked up a further bump plus flipped to sold out.
$ws1.Cells.Item(11, 6).Value = 28186
$ws1.Cells.Item(11, 7).Value = "已售罄"

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$sheet2Edits = @(
    @(6, 369),
    @(7, 783),
    @(11, 4234),
    @(16, 43),
    @(22, 4230),
    @(24, 2)
)
foreach ($edit in $sheet2Edits) {
    $row = $edit[0]
    $val = $edit[1]
    $ws2.Cells.Item($row, 6).Value = $val
}

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(4, 6).Value = 1140

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$sheet4Edits = @(
    @(4, 1140),
    @(8, 369),
    @(10, 783),
    @(11, 1138),
    @(12, 1217),
    @(13, 163),
    @(14, 529),
    @(15, 295),
    @(17, 50),
    @(18, 1250),
    @(26, 237),
    @(27, 43),
    @(28, 43),
    @(29, 456),
    @(30, 10),
    @(33, 309),
    @(34, 593),
    @(35, 264),
    @(36, 339),
    @(38, 34),
    @(39, 651),
    @(41, 194),
    @(42, 92),
    @(45, 70),
    @(47, 597),
    @(51, 2)
)
foreach ($edit in $sheet4Edits) {
    $row = $edit[0]
    $val = $edit[1]
    $ws4.Cells.Item($row, 6).Value = $val
}
